# Checkpoint edit: add "audience" column to the RibaReport table, and add a new
# "distribution" sheet with a RibaDistribution table.
#
# Column-width note: this host's ColumnWidth setter/getter round-trips through
# Excel's internal "character width" <-> stored-XML-width conversion, which
# adds a constant offset of 11/12 (~0.9167) character units. To land on an
# exact target XML width W we must set ColumnWidth = W - 11/12.
$WIDTH_OFFSET = 11 / 12

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) report sheet: insert "audience" as the 4th column of the RibaReport table
#    (between class_no and rs), shifting rs/target/asp one column to the right.
# ---------------------------------------------------------------------------
$reportWs = $wb.Worksheets.Item("report")
$reportTable = $reportWs.ListObjects.Item(1)

# Grow the table by one column at the end (this also extends the sheet range
# to A1:G2) - we'll relabel/reshuffle the header text next.
$reportTable.ListColumns.Add() | Out-Null

# Shift the existing rs/target/asp header text one column to the right, then
# write the new "audience" header into column D. Writing directly into the
# header cells keeps the ListObject's column metadata in sync with the text.
$aspText = $reportWs.Cells.Item(1, 6).Text
$targetText = $reportWs.Cells.Item(1, 5).Text
$rsText = $reportWs.Cells.Item(1, 4).Text

$reportWs.Cells.Item(1, 7).Value = $aspText
$reportWs.Cells.Item(1, 6).Value = $targetText
$reportWs.Cells.Item(1, 5).Value = $rsText
$reportWs.Cells.Item(1, 4).Value = "audience"

$reportTable.Name = "RibaReport"
$reportTable.TableStyle = "TableStyleLight9"

# Column widths: A-C unchanged (10, 14, 10); D is the new "audience" column
# (12); E/F/G carry forward the old rs/target/asp widths (6, 55, 10).
$reportWs.Columns(4).ColumnWidth = 12 - $WIDTH_OFFSET
$reportWs.Columns(5).ColumnWidth = 6 - $WIDTH_OFFSET
$reportWs.Columns(6).ColumnWidth = 55 - $WIDTH_OFFSET
$reportWs.Columns(7).ColumnWidth = 10 - $WIDTH_OFFSET

# ---------------------------------------------------------------------------
# 2) Add the new "distribution" sheet right after "report".
# ---------------------------------------------------------------------------
$distWs = $wb.Worksheets.Add($null, $reportWs)
$distWs.Name = "distribution"

$distHeaders = @("level", "grade", "class_no", "audience", "item_no", "count_a", "count_b", "count_total", "pct_a", "pct_b")
for ($i = 0; $i -lt $distHeaders.Length; $i++) {
    $distWs.Cells.Item(1, $i + 1).Value = $distHeaders[$i]
}

$distTable = $distWs.ListObjects.Add(1, $distWs.Range("A1:J2"), $null, 1)
$distTable.Name = "RibaDistribution"
$distTable.TableStyle = "TableStyleLight9"

# Freeze the header row and match the selection/pane layout used on "report".
$distWs.Activate()
$distWs.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
$distWs.Range("A1").Select()

# Outline defaults matching the "report" sheet (summaryBelow / summaryRight).
$distWs.Outline.SummaryRow = 1
$distWs.Outline.SummaryColumn = 1

# Column widths for the distribution sheet.
$distWidths = @(10, 14, 10, 12, 8, 10, 10, 12, 10, 10)
for ($i = 0; $i -lt $distWidths.Length; $i++) {
    $distWs.Columns($i + 1).ColumnWidth = $distWidths[$i] - $WIDTH_OFFSET
}

# Restore focus to the report sheet (matches the original file's active tab).
$reportWs.Activate()
